# Atualizar mês de Maio
# Updates the "Preço Investido" (average invested price) column with the
# newly recalculated values for May, and rotates the last 7 tickers
# (rows 26-32) so each keeps the freshly computed running-average price
# that now belongs to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: updated average invested price values ---
$ws.Range("C3").Value  = 73.66131578947369
$ws.Range("C5").Value  = 12.88067375886525
$ws.Range("C7").Value  = 47.89415584415585
$ws.Range("C8").Value  = 74.83722222222222
$ws.Range("C9").Value  = 110.8544444444444
$ws.Range("C10").Value = 9.422576687116566
$ws.Range("C11").Value = 19.29623376623377
$ws.Range("C13").Value = 24.54693548387097
$ws.Range("C14").Value = 7.709543568464731
$ws.Range("C15").Value = 217.563
$ws.Range("C16").Value = 150.7195555555556
$ws.Range("C17").Value = 160.1547619047619
$ws.Range("C18").Value = 103.6566129032258
$ws.Range("C19").Value = 38.10470588235294
$ws.Range("C20").Value = 22.19833333333333
$ws.Range("C25").Value = 32.17232558139535

# --- Rows 26-32: tickers rotate up by one (COCA34 moves to the bottom) ---
$ws.Range("B26").Value = "AAPL34"
$ws.Range("C26").Value = 51.18026315789474

$ws.Range("B27").Value = "C1BS34"
$ws.Range("C27").Value = 68.45352941176471

$ws.Range("B28").Value = "MGLU3"
$ws.Range("C28").Value = 2.14

$ws.Range("B29").Value = "CPFF11"
$ws.Range("C29").Value = 72.45333333333333

$ws.Range("B30").Value = "ITSA2F"
$ws.Range("C30").Value = 2.98

$ws.Range("B31").Value = "VALE3F"
$ws.Range("C31").Value = 54.62790697674419

$ws.Range("B32").Value = "COCA34"
$ws.Range("C32").Value = 65.58620689655173
